$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet (so it lands after
# "SearchItems"), then rename it and give it the same "adidas" value that
# already exists elsewhere in the workbook's shared strings table.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "SearchItem"
$newSheet.Range("A2").Value = "adidas"

# Make the newly added sheet the active one, with A2 selected (matching
# the selection left behind on the source sheet).
$newSheet.Activate()
$newSheet.Range("A2").Select() | Out-Null
